$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.976.24"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "3.397.85"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.29"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.82"
$ws.Range("E6").Value = "  -1.77%  "
$ws.Range("D7").Value = "3.399.25"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.477"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.396"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "3.975.99"
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.19"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "3.391.79"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "61.066.89"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.86"
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("E21").Value = "  -5.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.32"
$ws.Range("E22").Value = "  -5.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.560"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.24"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  -4.45%  "
$ws.Range("D27").Value = "3.531.95"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.39"
$ws.Range("E30").Value = "  -3.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.02"
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.16"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  -1.93%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.01"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "167.84"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Value = "3.428.34"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.01"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.90"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.68"
$ws.Range("E46").Value = "  -3.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").Value = "2.478.60"
$ws.Range("E48").Value = "  -5.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.83"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.01"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").Value = "  +1.05%  "
Write-Output "Applied cryptos update"
